$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new schedule rows (40 and 41), continuing the existing A/B/C pattern
# (date in A, date in B, period description in C). Column A continues the
# weekly-date sequence already present in rows 1-39; column B/C hold the
# new period's target date and title.

$ws.Range("A40").Value = "2025/12/26"
$ws.Range("A41").Value = "2026/1/2"

$ws.Range("B40").Value = "2026/2/20"
$ws.Range("B41").Value = "2026/2/27"

$ws.Range("C40").Value = "第90期 第二代星途"
$ws.Range("C41").Value = "第91期 秘寶 開放區域 清涼地 祕寶效果: 透過元素袋獲取星途解讀道具有5%(18%) 翻倍"

# Keep view/selection state consistent with the edit (scroll down a bit and
# select the next empty cell, as the original author did after typing).
$excel.ActiveWindow.ScrollRow = 30
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C42").Select()
